{"js": "// Update the date line and all \"dividend\u00f7divisor=quotient, remainder\"\n// answer cells in the table to the new values from the target revision.\n// Each old string is unique within the document, so a plain body.search\n// for the exact old text followed by a full-text replace is safe and\n// preserves the existing run formatting (font/size) on every match.\nconst pairs = [\n  [\"2024-10-20 Sunday\", \"2024-10-21 Monday\"],\n  [\"243\u00f74=60, 3\", \"779\u00f79=86, 5\"],\n  [\"361\u00f75=72, 1\", \"724\u00f78=90, 4\"],\n  [\"655\u00f73=218, 1\", \"245\u00f76=40, 5\"],\n  [\"971\u00f75=194, 1\", \"630\u00f77=90, 0\"],\n  [\"895\u00f79=99, 4\", \"553\u00f79=61, 4\"],\n  [\"905\u00f72=452, 1\", \"859\u00f72=429, 1\"],\n  [\"170\u00f72=85, 0\", \"135\u00f76=22, 3\"],\n  [\"665\u00f78=83, 1\", \"916\u00f79=101, 7\"],\n  [\"638\u00f77=91, 1\", \"403\u00f77=57, 4\"],\n  [\"764\u00f73=254, 2\", \"941\u00f78=117, 5\"],\n  [\"489\u00f76=81, 3\", \"756\u00f73=252, 0\"],\n  [\"772\u00f79=85, 7\", \"316\u00f77=45, 1\"],\n  [\"834\u00f74=208, 2\", \"224\u00f79=24, 8\"],\n  [\"386\u00f77=55, 1\", \"347\u00f76=57, 5\"],\n  [\"526\u00f79=58, 4\", \"827\u00f75=165, 2\"],\n  [\"981\u00f73=327, 0\", \"146\u00f74=36, 2\"],\n  [\"324\u00f72=162, 0\", \"640\u00f78=80, 0\"],\n  [\"238\u00f72=119, 0\", \"734\u00f73=244, 2\"],\n  [\"686\u00f74=171, 2\", \"577\u00f73=192, 1\"],\n  [\"535\u00f73=178, 1\", \"700\u00f78=87, 4\"],\n  [\"153\u00f79=17, 0\", \"640\u00f78=80, 0\"],\n  [\"661\u00f76=110, 1\", \"100\u00f73=33, 1\"],\n  [\"878\u00f72=439, 0\", \"901\u00f74=225, 1\"],\n  [\"828\u00f78=103, 4\", \"689\u00f72=344, 1\"],\n  [\"996\u00f72=498, 0\", \"106\u00f77=15, 1\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date line and all \"dividend\u00f7divisor=quotient, remainder\"\n# answer cells in the table to the new values from the target revision.\n# Each old string is unique in the document, so Find/Replace on the\n# exact old text (MatchCase, WholeWord off, not wildcard) safely swaps\n# in the new text while keeping the existing run formatting intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-10-20 Sunday\", \"2024-10-21 Monday\"),\n    @(\"243\u00f74=60, 3\", \"779\u00f79=86, 5\"),\n    @(\"361\u00f75=72, 1\", \"724\u00f78=90, 4\"),\n    @(\"655\u00f73=218, 1\", \"245\u00f76=40, 5\"),\n    @(\"971\u00f75=194, 1\", \"630\u00f77=90, 0\"),\n    @(\"895\u00f79=99, 4\", \"553\u00f79=61, 4\"),\n    @(\"905\u00f72=452, 1\", \"859\u00f72=429, 1\"),\n    @(\"170\u00f72=85, 0\", \"135\u00f76=22, 3\"),\n    @(\"665\u00f78=83, 1\", \"916\u00f79=101, 7\"),\n    @(\"638\u00f77=91, 1\", \"403\u00f77=57, 4\"),\n    @(\"764\u00f73=254, 2\", \"941\u00f78=117, 5\"),\n    @(\"489\u00f76=81, 3\", \"756\u00f73=252, 0\"),\n    @(\"772\u00f79=85, 7\", \"316\u00f77=45, 1\"),\n    @(\"834\u00f74=208, 2\", \"224\u00f79=24, 8\"),\n    @(\"386\u00f77=55, 1\", \"347\u00f76=57, 5\"),\n    @(\"526\u00f79=58, 4\", \"827\u00f75=165, 2\"),\n    @(\"981\u00f73=327, 0\", \"146\u00f74=36, 2\"),\n    @(\"324\u00f72=162, 0\", \"640\u00f78=80, 0\"),\n    @(\"238\u00f72=119, 0\", \"734\u00f73=244, 2\"),\n    @(\"686\u00f74=171, 2\", \"577\u00f73=192, 1\"),\n    @(\"535\u00f73=178, 1\", \"700\u00f78=87, 4\"),\n    @(\"153\u00f79=17, 0\", \"640\u00f78=80, 0\"),\n    @(\"661\u00f76=110, 1\", \"100\u00f73=33, 1\"),\n    @(\"878\u00f72=439, 0\", \"901\u00f74=225, 1\"),\n    @(\"828\u00f78=103, 4\", \"689\u00f72=344, 1\"),\n    @(\"996\u00f72=498, 0\", \"106\u00f77=15, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
